{"js": "// Add the author's name/e-mail as a brand-new first paragraph of the\n// document (\"Dok\u00fcmantasyona isim soyisim eklendi.\" / \"Name and surname\n// were added to the documentation.\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst firstParagraph = paragraphs.items[0];\nfirstParagraph.insertParagraph(\"\u015eilan EK\u0130N silanekinceng@gmail.com\", Word.InsertLocation.before);\nawait context.sync();\n", "ps1": "# Add the author's name/e-mail as a brand-new first paragraph of the\n# document (\"Dok\u00fcmantasyona isim soyisim eklendi.\" / \"Name and surname\n# were added to the documentation.\").\n$d = $word.ActiveDocument\n\n$firstPara = $d.Paragraphs(1)\n$newRange = $firstPara.Range.InsertParagraphBefore()\n$firstPara.Range.Text = \"\u015eilan EK\u0130N silanekinceng@gmail.com\"\n"}
